$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Nama Developer" column (D) for rows 3-6, which were previously empty.
$ws.Range("D3").Value = "Eka Yunia & Yulisha Dian"
$ws.Range("D4").Value = "Yulisha Dian"
$ws.Range("D5").Value = "Eka Yunia & Yulisha Dian"
$ws.Range("D6").Value = "Eka Yunia "

# Widen column D to fit the new content, and switch off auto best-fit sizing.
$ws.Columns.Item(4).ColumnWidth = 22.3

# Update the active cell selection to D7, matching the author's final cursor position.
$ws.Range("D7").Select()
